# Auto-applies cached market-data values (currentAveragePrice / Leve profit
# columns H:N) for the rows identified in the upstream data refresh diff.
# Values are taken verbatim from the authoritative post-edit OOXML so that
# re-applying this script reproduces the same cell contents and, where the
# diff added or removed a cell, the same sparse row shape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4899.1816
$ws.Range("I51").Value = 3399.5
$ws.Range("K51").Value = 3399.5
$ws.Range("M51").Value = -2915.5

$ws.Range("H116").Value = 9438.200000000001
$ws.Range("I116").Value = 6969
$ws.Range("K116").Value = 6969
$ws.Range("M116").Value = -3527

$ws.Range("H125").Value = 5269.8
$ws.Range("I125").Value = 6774.75
$ws.Range("J125").Value = 4266.5
$ws.Range("K125").Value = 60972.75
$ws.Range("L125").Value = 38398.5
$ws.Range("M125").Value = -58512.75
$ws.Range("N125").Value = -43318.5

$ws.Range("H135").Value = 19231524
$ws.Range("I135").Value = 843.7143
$ws.Range("J135").Value = 41667316
$ws.Range("K135").Value = 7593.428699999999
$ws.Range("L135").Value = 375005844
$ws.Range("M135").Value = -5058.428699999999
$ws.Range("N135").Value = -375010914

$ws.Range("H138").Value = 1753.8125
$ws.Range("I138").Value = 1025.6
$ws.Range("J138").Value = 2967.5
$ws.Range("K138").Value = 3076.8
$ws.Range("L138").Value = 8902.5
$ws.Range("M138").Value = 2063.2
$ws.Range("N138").Value = -19182.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3125.4546
$ws.Range("I2").Value = 1654.8572
$ws.Range("K2").Value = 1654.8572
$ws.Range("M2").Value = -1541.8572

$ws.Range("H32").Value = 4646.5557
$ws.Range("J32").Value = 9499.666999999999
$ws.Range("L32").Value = 9499.666999999999
$ws.Range("N32").Value = -10073.667

$ws.Range("H61").Value = 29413990
$ws.Range("I61").Value = 33335322
$ws.Range("K61").Value = 33335322
$ws.Range("M61").Value = -33335110

$ws.Range("H63").Value = 2050.3333
$ws.Range("I63").Value = 1760.4
$ws.Range("K63").Value = 1760.4
$ws.Range("M63").Value = -1074.4

$ws.Range("H66").Value = 2050.3333
$ws.Range("I66").Value = 1760.4
$ws.Range("K66").Value = 8802
$ws.Range("M66").Value = -5370

$ws.Range("H97").Value = 3217.9092
$ws.Range("I97").Value = 3233.0476
$ws.Range("J97").Value = 2900
$ws.Range("K97").Value = 3233.0476
$ws.Range("L97").Value = 2900
$ws.Range("M97").Value = -2737.0476
$ws.Range("N97").Value = -3892

$ws.Range("H110").Value = 799
$ws.Range("I110").Value = 799
$ws.Range("K110").Value = 799
$ws.Range("M110").Value = 1246

$ws.Range("H116").Value = 3125.4546
$ws.Range("I116").Value = 1654.8572
$ws.Range("K116").Value = 1654.8572
$ws.Range("M116").Value = 639.1428000000001

$ws.Range("H136").Value = 29413990
$ws.Range("I136").Value = 33335322
$ws.Range("K136").Value = 100005966
$ws.Range("M136").Value = -100003416

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3125.4546
$ws.Range("I3").Value = 1654.8572
$ws.Range("K3").Value = 1654.8572
$ws.Range("M3").Value = -1540.8572

$ws.Range("H99").Value = 4973
$ws.Range("I99").Value = 2916.3333
$ws.Range("K99").Value = 2916.3333
$ws.Range("M99").Value = -1418.3333

$ws.Range("H134").Value = 7000
$ws.Range("I134").Value = 7000
$ws.Range("K134").Value = 21000
$ws.Range("M134").Value = -18465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10003952
$ws.Range("I31").Value = 3156.2144
$ws.Range("K31").Value = 3156.2144
$ws.Range("M31").Value = -2861.2144

$ws.Range("H34").Value = 10003952
$ws.Range("I34").Value = 3156.2144
$ws.Range("K34").Value = 3156.2144
$ws.Range("M34").Value = -2954.2144

$ws.Range("H132").Value = 67594.77
$ws.Range("I132").Value = 71925.86
$ws.Range("J132").Value = 4794
$ws.Range("K132").Value = 215777.58
$ws.Range("L132").Value = 14382
$ws.Range("M132").Value = -213247.58
$ws.Range("N132").Value = -19442

$ws.Range("H134").Value = 2815.2334
$ws.Range("I134").Value = 2384.5217
$ws.Range("J134").Value = 4230.4287
$ws.Range("K134").Value = 7153.5651
$ws.Range("L134").Value = 12691.2861
$ws.Range("M134").Value = -4618.5651
$ws.Range("N134").Value = -17761.2861

$ws.Range("H141").Value = 101354.8
$ws.Range("J141").Value = 118568.5
$ws.Range("L141").Value = 118568.5
$ws.Range("N141").Value = -128928.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 485.75
$ws.Range("I26").Value = 261.8
$ws.Range("K26").Value = 785.4000000000001
$ws.Range("M26").Value = -497.4000000000001

$ws.Range("H32").Value = 1000000
$ws.Range("I32").Value = 1000000
$ws.Range("K32").Value = 3000000
$ws.Range("M32").Value = -2999717

$ws.Range("H39").Value = 6785.933
$ws.Range("I39").Value = 906.8461
$ws.Range("J39").Value = 45000
$ws.Range("K39").Value = 2720.5383
$ws.Range("L39").Value = 135000
$ws.Range("M39").Value = -2426.5383
$ws.Range("N39").Value = -135588

$ws.Range("H96").Value = 4174.3335
$ws.Range("I96").Value = 3762
$ws.Range("K96").Value = 11286
$ws.Range("M96").Value = -9227

$ws.Range("H101").Value = 14999
$ws.Range("J101").Value = 14999
$ws.Range("L101").Value = 44997
$ws.Range("N101").Value = -49865

$ws.Range("H138").Value = 6048.3335
$ws.Range("I138").Value = 6458
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 19374
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -14234
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1575.25
$ws.Range("I113").Value = 1575.25
$ws.Range("K113").Value = 1575.25
$ws.Range("M113").Value = 594.75

$ws.Range("H135").Value = 119994.5
$ws.Range("J135").Value = 119994.5
$ws.Range("L135").Value = 119994.5
$ws.Range("N135").Value = -130134.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4092.3333
$ws.Range("I100").Value = 3603.875
$ws.Range("K100").Value = 3603.875
$ws.Range("M100").Value = -3062.875

$ws.Range("H122").Value = 2980763.5
$ws.Range("I122").Value = 4003.2593
$ws.Range("J122").Value = 8338931.5
$ws.Range("K122").Value = 12009.7779
$ws.Range("L122").Value = 25016794.5
$ws.Range("M122").Value = -9559.777900000001
$ws.Range("N122").Value = -25021694.5

$ws.Range("H132").Value = 3433.861
$ws.Range("I132").Value = 3284.926
$ws.Range("K132").Value = 9854.778
$ws.Range("M132").Value = -7324.778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H113").Value = 845.2593000000001
$ws.Range("J113").Value = 756.36365
$ws.Range("L113").Value = 2269.09095
$ws.Range("N113").Value = -6609.09095

$ws.Range("H122").Value = 11114127
$ws.Range("I122").Value = 1853.7273
$ws.Range("J122").Value = 28576272
$ws.Range("K122").Value = 5561.1819
$ws.Range("L122").Value = 85728816
$ws.Range("M122").Value = -3111.1819
$ws.Range("N122").Value = -85733716

